# Minor edits on some file
# Slide 6, shape 13 ("Take things with a grain of salt ...") -
# rewrite the second run's wording and split it into three runs so the
# formatting-carrying run boundaries match the edited copy.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(6)
$sh = $s.Shapes.Item(13)
$tr = $sh.TextFrame.TextRange

# The paragraph currently holds two runs:
#   1) "Take things with a grain of salt"                              (32 chars)
#   2) " (somehow) as across all Organization, 69% of them have half
#       of their employee population took the survey which means it
#       is not indicative of the voice of the entire workforce"        (175 chars)
#
# Replace run 2's text with the corrected wording, then split it into
# three runs (matching the target edit) by re-setting each segment's
# text through TextRange.Characters(start, length).

$run2 = $tr.Characters(33, $tr.Length - 32)
$run2.Text = " as across all Organization, 69% of the employee population took the survey which means it is not indicative of the voice of the entire workforce"

$partB = $tr.Characters(73, 20)
$partB.Text = "employee population "

$partC = $tr.Characters(93, 85)
$partC.Text = "took the survey which means it is not indicative of the voice of the entire workforce"
